$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (hk_ump_polymer_stock) ---
$ws.Range("Q3").Value = 0.21
$ws.Range("R3").Value = -9
$ws.Range("S3").Value = -13

# --- Row 4 (tommybuilt_ump_tailhook_brace_adapter) ---
$ws.Range("P4").Value = 15
$ws.Range("Q4").Value = 0.2
$ws.Range("R4").Value = -10
$ws.Range("S4").Value = -7

# --- Row 6 (ghw_tailhook_mod1_brace_small) ---
$ws.Range("P6").Value = -3
$ws.Range("Q6").Value = 0.08
$ws.Range("R6").Value = 2
$ws.Range("S6").Value = 2

# --- Row 7: magpul_moe_carbine_stock -> renamed to just "MOE Carbine", values punished ---
$ws.Range("A7").ClearContents()
$ws.Range("B7").Value = "MOE Carbine"

$ws.Range("D7").Value = 0.36
$ws.Range("E7").Value = -11
$ws.Range("F7").Value = -10
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 1400
$ws.Range("N7").Value = 24.6

$ws.Range("Q7").Value = 0.36
$ws.Range("R7").Value = -11
$ws.Range("S7").Value = -10
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 0
$ws.Range("X7").Value = 0
$ws.Range("Y7").Value = 0
$ws.Range("Z7").Value = 1400
$ws.Range("AA7").Value = 24.6

# --- Row 9 (tommybuilt_ump_tailhook_brace_adapter dup/other) ---
$ws.Range("Q9").Value = 0.13

# --- Row 10 (ghw_tailhook_mod1_brace_small dup/other) ---
$ws.Range("P10").Value = 12

# --- Sheet selection ---
$ws.Range("H18").Select()
